# Move the "Comment" column (J) to the end of the block (after column O),
# shifting WaitingTime/NumberMeasuringFields/L/A/B (K:O) left into J:N.
# Applies to the header row and the 4 metadata rows below it (rows 1-5).
#
# Row 1 (headers) needs a full 6-way circular shift since every label differs.
# Rows 2-5 only ever have real text in column J (the rest of the block, K:O,
# is blank) so the whole edit there reduces to swapping J and O - the blank
# K:N cells are left completely untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("J", "K", "L", "M", "N", "O")

# Row 1: full circular shift left.
$vals = @{}
foreach ($c in $cols) {
    $vals[$c] = $ws.Range("${c}1").Value()
}
$ws.Range("J1").Value = $vals["K"]
$ws.Range("K1").Value = $vals["L"]
$ws.Range("L1").Value = $vals["M"]
$ws.Range("M1").Value = $vals["N"]
$ws.Range("N1").Value = $vals["O"]
$ws.Range("O1").Value = $vals["J"]

# Rows 2-5: only J (comment text) and O (previously blank) actually change.
for ($r = 2; $r -le 5; $r++) {
    $jVal = $ws.Range("J$r").Value()
    $oVal = $ws.Range("O$r").Value()
    if ("$jVal" -ne "$oVal") {
        $ws.Range("O$r").Value = $jVal
        if ("$jVal" -ne "") {
            $ws.Range("J$r").ClearContents()
        }
    }
}
